$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the values of columns D, N, O, P, Q, R, S, T between
# row 2 and row 3 (columns A, B, C, E-M are identical between the two
# rows, so only these columns actually change).

$cols = @("D", "N", "O", "P", "Q", "R", "S", "T")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $val3
    $ws.Range($addr3).Value = $val2
}
